$d = $word.ActiveDocument

# The "Don vi thi cong" acceptance table (3rd column holds the {#job}... text)
# is the 4th table in the document. Resize it from an auto-fit table to a
# fixed-width table (8651 dxa = 432.55 pt) and widen its 3rd column from
# 3645 dxa (182.25 pt) to 3933 dxa (196.65 pt).
$t = $d.Tables.Item(4)

$t.PreferredWidthType = 3   # wdPreferredWidthPoints
$t.PreferredWidth = 432.55  # 8651 dxa

$t.Columns.Item(3).Width = 196.65  # 3933 dxa
